# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Total "VALOR MORA" header figure
$ws.Range("E11").Value = 284700

# Shift the 5 period columns forward by one month: drop 2503, add 2508.
# (2507 -> row16, 2506 -> row17 ... become 2504,2505,2506,2507,2508)
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2508"

# Valor Mora column (G) updated for all 5 rows
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
$ws.Range("G19").Value = 1423500
$ws.Range("G20").Value = 1423500

# Salario Basico column (F) for the newly added last period row
$ws.Range("F20").Value = 56940
